$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.75
$ws.Range("K3").Value = 1.83
$ws.Range("L3").Value = 5.5
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25
$ws.Range("Q3").Value = 2.88
$ws.Range("R3").Value = 1.4
$ws.Range("S3").Value = 1.67
$ws.Range("T3").Value = 2.1
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 4.75
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 5.5
$ws.Range("AE3").Value = 23
$ws.Range("AF3").Value = 101
$ws.Range("AI3").Value = 17
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 51
$ws.Range("AP3").Value = 29
$ws.Range("AR3").Value = 81
$ws.Range("AS3").Value = 351
$ws.Range("AT3").Value = 2.1
$ws.Range("AU3").Value = 10
$ws.Range("AX3").Value = 29
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 126
$ws.Range("G10").Value = 2.45
$ws.Range("J10").Value = 3.25
$ws.Range("M10").Value = 1.11
$ws.Range("N10").Value = 6.5
$ws.Range("U10").Value = 2.1
$ws.Range("V10").Value = 1.67
$ws.Range("AG10").Value = 7.5
$ws.Range("AW10").Value = 5
$ws.Range("I14").Value = 1.57
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 9
$ws.Range("Z14").Value = 67
$ws.Range("AC14").Value = 9
$ws.Range("Q15").Value = 2.4
$ws.Range("R15").Value = 1.53
$ws.Range("G17").Value = 1.85
$ws.Range("I17").Value = 3.8
$ws.Range("J17").Value = 2.6
$ws.Range("L17").Value = 4.5
$ws.Range("N17").Value = 10
$ws.Range("Z17").Value = 15
$ws.Range("AA17").Value = 15
$ws.Range("AH17").Value = 21
$ws.Range("AO17").Value = 10
$ws.Range("AW17").Value = 6
$ws.Range("AZ17").Value = 81
$ws.Range("G18").Value = 7.9
$ws.Range("H18").Value = 5.1
$ws.Range("I18").Value = 1.32
$ws.Range("K18").Value = 2.77
$ws.Range("L18").Value = 1.7
$ws.Range("O18").Value = 1.09
$ws.Range("P18").Value = 6.2
$ws.Range("Q18").Value = 1.3
$ws.Range("R18").Value = 3.25
$ws.Range("S18").Value = 1.19
$ws.Range("T18").Value = 4.2
$ws.Range("U18").Value = 1.5
$ws.Range("V18").Value = 2.42
$ws.Range("W18").Value = 40
$ws.Range("X18").Value = 80
$ws.Range("Y18").Value = 25
$ws.Range("Z18").Value = 200
$ws.Range("AA18").Value = 70
$ws.Range("AB18").Value = 40
$ws.Range("AC18").Value = 10.75
$ws.Range("AD18").Value = 12
$ws.Range("AE18").Value = 15
$ws.Range("AF18").Value = 37
$ws.Range("AG18").Value = 13
$ws.Range("AH18").Value = 9.75
$ws.Range("AI18").Value = 9
$ws.Range("AJ18").Value = 10.75
$ws.Range("AK18").Value = 9.75
$ws.Range("AL18").Value = 16.5
$ws.Range("AM18").Value = 175
$ws.Range("AN18").Value = 10
$ws.Range("AP18").Value = 26
$ws.Range("AR18").Value = 150
$ws.Range("AS18").Value = 175
$ws.Range("AT18").Value = 4.2
$ws.Range("AU18").Value = 6.9
$ws.Range("AV18").Value = 35
$ws.Range("AW18").Value = 3.8
$ws.Range("AX18").Value = 6
$ws.Range("AY18").Value = 11
$ws.Range("AZ18").Value = 13.5
$ws.Range("BA18").Value = 25
$ws.Range("BB18").Value = 80
$ws.Range("BC18").Value = 400
